$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLR")

# Insert a new column before column D. This shifts the existing D:K data
# (and formatting) right by one column to E:L, exactly like the in-app
# Insert Column command.
$ws.Columns("D").Insert()

# The newly inserted column D has no formatting yet; clone number/date
# formats from column E (which now holds what used to be column D) so the
# new column renders consistently with the rest of the table.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period's figures.
$newPeriodValues = @(
    @{Row=7; Value=43465},
    @{Row=8; Value=19166600},
    @{Row=9; Value=18496700},
    @{Row=10; Value=669900},
    @{Row=13; Value=0},
    @{Row=15; Value=0},
    @{Row=17; Value=18684800},
    @{Row=18; Value=481800},
    @{Row=20; Value=0},
    @{Row=21; Value=698400},
    @{Row=22; Value=0},
    @{Row=23; Value=481800},
    @{Row=24; Value=188800},
    @{Row=25; Value=0},
    @{Row=26; Value=293000},
    @{Row=27; Value=224800},
    @{Row=28; Value=0},
    @{Row=30; Value=0},
    @{Row=31; Value=0},
    @{Row=32; Value=0},
    @{Row=33; Value=224800},
    @{Row=34; Value=0},
    @{Row=35; Value=224800},
    @{Row=38; Value=43465},
    @{Row=41; Value=1764700},
    @{Row=42; Value=214800},
    @{Row=43; Value=1534300},
    @{Row=44; Value=1545000},
    @{Row=45; Value=382000},
    @{Row=46; Value=5440900},
    @{Row=47; Value=938500},
    @{Row=48; Value=1013700},
    @{Row=49; Value=533600},
    @{Row=50; Value=0},
    @{Row=51; Value=0},
    @{Row=52; Value=986900},
    @{Row=53; Value=0},
    @{Row=54; Value=8913600},
    @{Row=57; Value=1638900},
    @{Row=58; Value=26900},
    @{Row=59; Value=1886700},
    @{Row=60; Value=3552500},
    @{Row=61; Value=1661600},
    @{Row=62; Value=581500},
    @{Row=63; Value=0},
    @{Row=64; Value=0},
    @{Row=65; Value=0},
    @{Row=66; Value=5950500},
    @{Row=68; Value=0},
    @{Row=69; Value=0},
    @{Row=70; Value=0},
    @{Row=71; Value=0},
    @{Row=72; Value=3422200},
    @{Row=73; Value=0},
    @{Row=74; Value=0},
    @{Row=75; Value=0},
    @{Row=76; Value=2963200},
    @{Row=77; Value=0},
    @{Row=80; Value=43465},
    @{Row=81; Value=224800},
    @{Row=83; Value=216700},
    @{Row=84; Value=0},
    @{Row=85; Value=0},
    @{Row=86; Value=0},
    @{Row=87; Value=0},
    @{Row=88; Value=0},
    @{Row=89; Value=162200},
    @{Row=91; Value=-211000},
    @{Row=92; Value=0},
    @{Row=93; Value=0},
    @{Row=94; Value=1400},
    @{Row=96; Value=-118700},
    @{Row=97; Value=0},
    @{Row=98; Value=0},
    @{Row=99; Value=0},
    @{Row=100; Value=-140500},
    @{Row=101; Value=-62400},
    @{Row=102; Value=-39300}
)
foreach ($item in $newPeriodValues) {
    $ws.Cells.Item($item.Row, 4).Value = $item.Value
}

# Rows where the latest period has no applicable figure ("NA").
$naRows = @(12, 14, 29)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 4).Value = "NA"
}
